# Acid_test.xlsx update: fill in the SMILES column (E) for each acid row
# and update the last active cell selection, matching the author's edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("test")

$ws.Range("E2").Value = "Cl"
$ws.Range("E3").Value = "OS(=O)(=O)O"
$ws.Range("E4").Value = "ON(=O)=O"
$ws.Range("E5").Value = "OC(=O)CC(O)(C(=O)O)CC(=O)O"
$ws.Range("E6").Value = "CC(O)=O"

$ws.Activate()
$ws.Range("G11").Select()
